# Append new lancers.jp listings scraped at 2025-11-05 18:25:40 (JST).
# Two brand-new rows are inserted right before the former row 9, pushing
# the former rows 9-10 down to 11-12. The "取得日時" (fetched-at) timestamp
# on every existing data row is refreshed to the new run's timestamp too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-05 18:25:40"

# --- 1. Refresh the timestamp column for the rows 2-8, which keep their
#        position and all other values unchanged ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 2. Shift the old rows 9-10 down to 11-12 (literal values, since the
#        former contents of those rows are already known from the sheet) ---
$ws.Range("A11").Value = $newTimestamp
$ws.Range("B11").Value = "【継続依頼あり】GASやn8nのオンラインセミナー研修講師を募集!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5427459"
$ws.Range("G11").Value = 13

$ws.Range("A12").Value = $newTimestamp
$ws.Range("B12").Value = "中国語 ワードプレスの分かる人"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "~ 5,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5427699"
$ws.Range("G12").Value = 10

# --- 3. Write the two brand-new rows into the now-vacated 9 and 10 slots ---
$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5427682"
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = "◆開発"

$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = "【急募】完全オンラインでのLstep構築・運用依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5427793"
$ws.Range("G10").Value = 18

# H10, H11, H12 stay blank, matching the source listings that had no
# "skill summary" tag.

# --- 4. Rebuild the URL hyperlinks for column F (rows 2-12) from scratch.
#        Row insertion/shifting in this engine does not re-anchor the
#        Hyperlinks collection, so clear everything and re-add in order,
#        using the same literal URLs that were just written into column F. ---
$ws.Cells.Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5416301"
    3  = "https://www.lancers.jp/work/detail/5420440"
    4  = "https://www.lancers.jp/work/detail/5416328"
    5  = "https://www.lancers.jp/work/detail/5427011"
    6  = "https://www.lancers.jp/work/detail/5427648"
    7  = "https://www.lancers.jp/work/detail/5427397"
    8  = "https://www.lancers.jp/work/detail/5427338"
    9  = "https://www.lancers.jp/work/detail/5427682"
    10 = "https://www.lancers.jp/work/detail/5427793"
    11 = "https://www.lancers.jp/work/detail/5427459"
    12 = "https://www.lancers.jp/work/detail/5427699"
}

foreach ($r in 2,3,4,5,6,7,8,9,10,11,12) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
}
